$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a serial date value that was bumped
# by one day (46060 -> 46061) for every data row (rows 2-505).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 505 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 46061
